# Update "horarios" (schedule) workbook with the latest scrape results.
# - Línea 141 (LP1912) sheet: refresh timestamp/total and rewrite the schedule rows (13 -> 17 rows).
# - Línea 141 (LP1912-215) sheet: refresh timestamp/total and rewrite the schedule rows (2 -> 3 rows).
# - Línea 141 (6203-6173) sheet: refresh timestamp only (still 0 rows).

$wb = $excel.ActiveWorkbook

$newTimestamp = "05:42:52"

# -----------------------------------------------------------------
# Sheet "LP1912"
# -----------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("LP1912")

$ws1.Range("A2").Value = "Última actualización: " + $newTimestamp
$ws1.Range("A3").Value = "Total filas: 17"

$data1 = @(
    @("05:42:52", "05:43", "14_ABASTO",     1,   "LP1912"),
    @("05:42:52", "05:46", "17_ROMERO",     4,   "LP1912"),
    @("05:42:52", "06:08", "10_OLMOS",      26,  "LP1912"),
    @("05:42:52", "06:15", "215A_EL PATO",  33,  "LP1912"),
    @("05:42:52", "06:29", "23_HERNANDEZ",  47,  "LP1912"),
    @("05:42:52", "06:33", "11_ETCHEVERRY", 51,  "LP1912"),
    @("05:42:52", "06:38", "17X38_ROMERO",  56,  "LP1912"),
    @("05:42:52", "06:40", "16_SANTA ANA",  58,  "LP1912"),
    @("05:42:52", "06:56", "215A_EL PATO",  74,  "LP1912"),
    @("05:42:52", "06:58", "225_GOMEZ",     76,  "LP1912"),
    @("05:42:52", "07:15", "215C_EL PATO",  93,  "LP1912"),
    @("05:42:52", "07:18", "14_ABASTO",     96,  "LP1912"),
    @("05:42:52", "07:20", "16_SANTA ANA",  98,  "LP1912"),
    @("05:42:52", "07:21", "23_HERNANDEZ",  99,  "LP1912"),
    @("05:42:52", "07:29", "17X38_ROMERO",  107, "LP1912"),
    @("05:42:52", "07:34", "10_OLMOS",      112, "LP1912"),
    @("05:42:52", "07:36", "27_EL RETIRO",  114, "LP1912")
)

$row = 6
foreach ($r in $data1) {
    $ws1.Cells.Item($row, 1).Value = $r[0]
    $ws1.Cells.Item($row, 2).Value = $r[1]
    $ws1.Cells.Item($row, 3).Value = $r[2]
    $ws1.Cells.Item($row, 4).Value = $r[3]
    $ws1.Cells.Item($row, 5).Value = $r[4]
    $row = $row + 1
}

# -----------------------------------------------------------------
# Sheet "LP1912-215"
# -----------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")

$ws2.Range("A2").Value = "Última actualización: " + $newTimestamp
$ws2.Range("A3").Value = "Total filas: 3"

$data2 = @(
    @("05:42:52", "06:15", "215A_EL PATO", 33, "LP1912"),
    @("05:42:52", "06:56", "215A_EL PATO", 74, "LP1912"),
    @("05:42:52", "07:15", "215C_EL PATO", 93, "LP1912")
)

$row = 6
foreach ($r in $data2) {
    $ws2.Cells.Item($row, 1).Value = $r[0]
    $ws2.Cells.Item($row, 2).Value = $r[1]
    $ws2.Cells.Item($row, 3).Value = $r[2]
    $ws2.Cells.Item($row, 4).Value = $r[3]
    $ws2.Cells.Item($row, 5).Value = $r[4]
    $row = $row + 1
}

# -----------------------------------------------------------------
# Sheet "6203-6173" (only the timestamp refreshes, still 0 data rows)
# -----------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("6203-6173")

$ws3.Range("A2").Value = "Última actualización: " + $newTimestamp

Write-Host "Horarios actualizados Linea 141"
